$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.967.20'
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.821.26'
$ws.Range("E3").Value = '  +2.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '355.11'
$ws.Range("E5").Value = '  +6.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '113.88'
$ws.Range("E6").Value = '  -2.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.553'
$ws.Range("E7").Value = '  +2.62%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.599'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.61'
$ws.Range("E10").Value = '  -0.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0854'
$ws.Range("E11").Value = '  -0.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.01'
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("E13").Value = '  +1.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.74'
$ws.Range("E14").Value = '  +1.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.260.89'
$ws.Range("E15").Value = '  +2.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.827.68'
$ws.Range("E16").Value = '  +2.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.896'
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.838.87'
$ws.Range("E18").Value = '  +0.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.39'
$ws.Range("E19").Value = '  +7.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.16'
$ws.Range("E20").Value = '  -1.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.54'
$ws.Range("E21").Value = '  +0.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0999'
$ws.Range("E22").Value = '  +2.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '270.30'
$ws.Range("E23").Value = '  -2.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.85'
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.79'
$ws.Range("E25").Value = '  +5.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.76'
$ws.Range("E26").Value = '  -0.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.31'
$ws.Range("E28").Value = '  +1.07%  '
$ws.Range("E29").Value = '  +1.83%  '
$ws.Range("E30").Value = '  -0.25%  '
$ws.Range("B31").Value = 'VeChain'
$ws.Range("C31").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0456'
$ws.Range("E31").Value = '  +30.13%  '
$ws.Range("B32").Value = 'OKB'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.65'
$ws.Range("E32").Value = '  +1.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '33.90'
$ws.Range("E33").Value = '  -3.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0832'
$ws.Range("E35").Value = '  +0.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  -0.09%  '
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.22'
$ws.Range("E38").Value = '  -0.81%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.87'
$ws.Range("E39").Value = '  -2.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.20'
$ws.Range("E40").Value = '  -4.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.72'
$ws.Range("E41").Value = '  +2.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.57'
$ws.Range("E42").Value = '  +5.22%  '
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.115'
$ws.Range("E43").Value = '  +1.18%  '
$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '126.13'
$ws.Range("E44").Value = '  -0.85%  '
$ws.Range("E45").Value = '  +0.33%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.36'
$ws.Range("E46").Value = '  +1.11%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.084.91'
$ws.Range("E47").Value = '  -0.27%  '
$ws.Range("E48").Value = '  +3.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.68'
$ws.Range("E49").Value = '  +2.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.948'
$ws.Range("E50").Value = '  +8.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '60.66'
$ws.Range("E51").Value = '  +1.21%  '
